$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos in the "topic" column (F) by re-typing the cell text. Excel will
# append fresh shared-string entries and drop the old mis-spelled ones once
# they are no longer referenced anywhere in the workbook.
$ws.Range("F6").Value = "Integer Multiplication Division"
$ws.Range("F2").Value = "Integer Multiplication Max Min -1"
$ws.Range("F3").Value = "Integer Multiplication Max Min -2"

# Remove the now-unused trailing blank rows (11-19) below the data table.
$ws.Range("A11:K19").EntireRow.Delete()

# Leave the cursor where the author's last edit landed.
$ws.Range("F18").Select()
